$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Excel table row number" column moves from B (with the old
# literal row index kept in A) to A, and the whole table shifts one
# column to the left (B..G -> A..F). Column G is dropped entirely and
# 6 more data rows are appended (rows 14-19).
# ------------------------------------------------------------------

# Grab the existing header format (currently applied to B1:G1, style
# index 1 - bold, bordered, centered) and stamp it onto A1 before we
# touch any of the existing cells, so A1 ends up re-using that exact
# style.
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Header row text, now living in A1:F1.
$ws.Range("A1").Value = "Excel table row number"
$ws.Range("B1").Value = "Question name"
$ws.Range("C1").Value = "Crossbreak subgroup"
$ws.Range("D1").Value = "National average for question (%)"
$ws.Range("E1").Value = "Proportion for subgroup (%)"
$ws.Range("F1").Value = "Significant difference (%)"

# Old column G is no longer part of the table - drop its contents and
# formatting completely (header + all former data rows).
$ws.Range("G1:G13").Clear()

# The old column A (rows 2-13) held a plain row index with the bold
# bordered header style accidentally applied to it; the new layout has
# no styling on the data cells at all, so strip that formatting.
$ws.Range("A2:A13").ClearFormats()

# Now populate the full corrected data block, A2:F19.
$ws.Range("A2").Value = 97
$ws.Range("B2").Value = "Retired before the pandemic"
$ws.Range("C2").Value = "18-24"
$ws.Range("D2").Value = 43.770211919328
$ws.Range("E2").Value = 0.367092365260114
$ws.Range("F2").Value = -43.40311955406789

$ws.Range("A3").Value = 97
$ws.Range("B3").Value = "Retired before the pandemic"
$ws.Range("C3").Value = "45-54"
$ws.Range("D3").Value = 43.770211919328
$ws.Range("E3").Value = 3.49356887781973
$ws.Range("F3").Value = -40.27664304150827

$ws.Range("A4").Value = 99
$ws.Range("B4").Value = "Currently focusing on education or study"
$ws.Range("C4").Value = "18-24"
$ws.Range("D4").Value = 9.356530781419
$ws.Range("E4").Value = 52.3122273459771
$ws.Range("F4").Value = 42.9556965645581

$ws.Range("A5").Value = 111
$ws.Range("B5").Value = "Very likely"
$ws.Range("C5").Value = "18-24"
$ws.Range("D5").Value = 20.0082371014857
$ws.Range("E5").Value = 77.76134286736071
$ws.Range("F5").Value = 57.753105765875

$ws.Range("A6").Value = 115
$ws.Range("B6").Value = "Very unlikely"
$ws.Range("C6").Value = "18-24"
$ws.Range("D6").Value = 58.86715293352231
$ws.Range("E6").Value = 4.24767106839781
$ws.Range("F6").Value = -54.6194818651245

$ws.Range("A7").Value = 115
$ws.Range("B7").Value = "Very unlikely"
$ws.Range("C7").Value = "25-34"
$ws.Range("D7").Value = 58.86715293352231
$ws.Range("E7").Value = 7.84713647937069
$ws.Range("F7").Value = -51.02001645415162

$ws.Range("A8").Value = 117
$ws.Range("B8").Value = "Total Likely:"
$ws.Range("C8").Value = "18-24"
$ws.Range("D8").Value = 26.477548712461
$ws.Range("E8").Value = 88.73009773374551
$ws.Range("F8").Value = 62.25254902128449

$ws.Range("A9").Value = 118
$ws.Range("B9").Value = "Total Unlikely:"
$ws.Range("C9").Value = "18-24"
$ws.Range("D9").Value = 64.4043291531944
$ws.Range("E9").Value = 5.18587048384209
$ws.Range("F9").Value = -59.2184586693523

$ws.Range("A10").Value = 118
$ws.Range("B10").Value = "Total Unlikely:"
$ws.Range("C10").Value = "25-34"
$ws.Range("D10").Value = 64.4043291531944
$ws.Range("E10").Value = 10.0463291479188
$ws.Range("F10").Value = -54.3580000052756

$ws.Range("A11").Value = 119
$ws.Range("B11").Value = "Net:"
$ws.Range("C11").Value = "18-24"
$ws.Range("D11").Value = -37.9267804407334
$ws.Range("E11").Value = 83.5442272499034
$ws.Range("F11").Value = 121.4710076906368

$ws.Range("A12").Value = 119
$ws.Range("B12").Value = "Net:"
$ws.Range("C12").Value = "25-34"
$ws.Range("D12").Value = -37.9267804407334
$ws.Range("E12").Value = 55.99427110512691
$ws.Range("F12").Value = 93.92105154586031

$ws.Range("A13").Value = 119
$ws.Range("B13").Value = "Net:"
$ws.Range("C13").Value = "35-44"
$ws.Range("D13").Value = -37.9267804407334
$ws.Range("E13").Value = 21.0616633985113
$ws.Range("F13").Value = 58.9884438392447

$ws.Range("A14").Value = 119
$ws.Range("B14").Value = "Net:"
$ws.Range("C14").Value = "65+"
$ws.Range("D14").Value = -37.9267804407334
$ws.Range("E14").Value = -94.85353342083511
$ws.Range("F14").Value = -56.9267529801017

$ws.Range("A15").Value = 119
$ws.Range("B15").Value = "Net:"
$ws.Range("C15").Value = "London"
$ws.Range("D15").Value = -37.9267804407334
$ws.Range("E15").Value = 8.31705680363695
$ws.Range("F15").Value = 46.24383724437035

$ws.Range("A16").Value = 1085
$ws.Range("B16").Value = "Somewhat more difficult"
$ws.Range("C16").Value = "Wales"
$ws.Range("D16").Value = 25.7372170148127
$ws.Range("E16").Value = 77.6090706348227
$ws.Range("F16").Value = 51.87185362000999

$ws.Range("A17").Value = 1086
$ws.Range("B17").Value = "Would not have made much difference"
$ws.Range("C17").Value = "35-44"
$ws.Range("D17").Value = 8.169964514875131
$ws.Range("E17").Value = 49.505442985824
$ws.Range("F17").Value = 41.33547847094887

$ws.Range("A18").Value = 1086
$ws.Range("B18").Value = "Would not have made much difference"
$ws.Range("C18").Value = "45-54"
$ws.Range("D18").Value = 8.169964514875131
$ws.Range("E18").Value = 59.1162973399879
$ws.Range("F18").Value = 50.94633282511277

$ws.Range("A19").Value = 1100
$ws.Range("B19").Value = "Significantly more difficult"
$ws.Range("C19").Value = "East of England"
$ws.Range("D19").Value = 35.4535835172556
$ws.Range("E19").Value = 85.3524049118576
$ws.Range("F19").Value = 49.898821394602
